$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (Serie date, FPD, FPL) -- "Actualización desde MV -datos-"
$newRows = @(
    @("08-10-2021", 4196150, 439805),
    @("12-10-2021", 844150, 1173581),
    @("13-10-2021", 956470, 1411844),
    @("14-10-2021", 5144270, 175237)
)

# Find first empty row after existing data (row 195 in this workbook)
$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i

    # Column A holds a "dd-mm-yyyy" style label. Enter it as an explicit text
    # formula first (so values like "08-10-2021" are not auto-parsed as a
    # date by Excel's input heuristics), then convert it to a plain stored
    # value via copy / paste-values so the cell keeps its default styling.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Formula = '="' + $newRows[$i][0] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}

$excel.CutCopyMode = 0
